$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "mapsto" column (C) with the new descriptive stage names
# (previously these were raw numeric codes).
$ws.Range("C2").Value = "wake"
$ws.Range("C3").Value = "rem"
$ws.Range("C4").Value = "stage1"
$ws.Range("C5").Value = "stage2"
$ws.Range("C6").Value = "sws"
$ws.Range("C7").Value = "sws"
$ws.Range("C8").Value = "artifact"
$ws.Range("C9").Value = "unknown"

# Update the selected cell to match the authored state
$ws.Range("C10").Select()
